$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force D:E columns to text first so numeric-looking values (e.g. "207.24")
# are stored as text, matching the original inlineStr cell type.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "27.176.05"
$ws.Range("E2").Value = "  -0.96%  "

$ws.Range("D3").Value = "1.574.17"
$ws.Range("E3").Value = "  -0.25%  "

$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").Value = "207.24"
$ws.Range("E5").Value = "  -0.19%  "

$ws.Range("E6").Value = "  -1.34%  "

$ws.Range("E7").Value = "  -0.13%  "

$ws.Range("D8").Value = "22.28"
$ws.Range("E8").Value = "  +0.15%  "

$ws.Range("E9").Value = "  -0.59%  "

$ws.Range("E10").Value = "  -0.01%  "

$ws.Range("E11").Value = "  +0.37%  "

$ws.Range("D12").Value = "1.796.85"
$ws.Range("E12").Value = "  -0.37%  "

$ws.Range("D13").Value = "1.578.48"
$ws.Range("E13").Value = "  +0.07%  "

$ws.Range("E14").Value = "  -1.23%  "

$ws.Range("E15").Value = "  -0.98%  "

$ws.Range("D16").Value = "27.181.06"
$ws.Range("E16").Value = "  -1.01%  "

$ws.Range("D17").Value = "62.28"
$ws.Range("E17").Value = "  -1.18%  "

$ws.Range("B18").Value = "Chainlink"
$ws.Range("C18").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D18").Value = "7.36"
$ws.Range("E18").Value = "  +0.98%  "

$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D19").Value = "214.11"
$ws.Range("E19").Value = "  -0.10%  "

$ws.Range("E20").Value = "  -0.71%  "

$ws.Range("E21").Value = "  -0.02%  "

$ws.Range("D22").Value = "4.13"
$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("D23").Value = "9.46"
$ws.Range("E23").Value = "  -3.16%  "

$ws.Range("D24").Value = "2.03"
$ws.Range("E24").Value = "  +1.27%  "

$ws.Range("D25").Value = "152.53"
$ws.Range("E25").Value = "  -0.47%  "

$ws.Range("D26").Value = "6.69"
$ws.Range("E26").Value = "  -2.99%  "

$ws.Range("D27").Value = "14.96"
$ws.Range("E27").Value = "  -0.66%  "

$ws.Range("E28").Value = "  -0.08%  "

$ws.Range("E29").Value = "  -1.16%  "

$ws.Range("E30").Value = "  -3.33%  "

$ws.Range("E31").Value = "  -1.66%  "

$ws.Range("E32").Value = "  -1.39%  "

$ws.Range("D33").Value = "1.404.26"
$ws.Range("E33").Value = "  +2.85%  "

$ws.Range("E34").Value = "  -1.18%  "

$ws.Range("D35").Value = "1.56"
$ws.Range("E35").Value = "  +2.04%  "

$ws.Range("E36").Value = "  -1.04%  "

$ws.Range("D37").Value = "0.942"
$ws.Range("E37").Value = "  -2.69%  "

$ws.Range("E38").Value = "  -1.65%  "

$ws.Range("D39").Value = "0.816"
$ws.Range("E39").Value = "  -0.49%  "

$ws.Range("E40").Value = "  -2.93%  "

$ws.Range("E41").Value = "  -0.04%  "

$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +3.00%  "

$ws.Range("D43").Value = "1.85"
$ws.Range("E43").Value = "  +4.07%  "

$ws.Range("D44").Value = "5.43"
$ws.Range("E44").Value = "  +3.19%  "

$ws.Range("E45").Value = "  +1.16%  "

$ws.Range("D46").Value = "63.74"
$ws.Range("E46").Value = "  -0.66%  "

$ws.Range("D47").Value = "1.709.39"
$ws.Range("E47").Value = "  -0.40%  "

$ws.Range("D48").Value = "85.65"
$ws.Range("E48").Value = "  -0.39%  "

$ws.Range("D49").Value = "0.0₇0991"
$ws.Range("E49").Value = "  -0.21%  "

$ws.Range("D50").Value = "0.0952"
$ws.Range("E50").Value = "  -0.50%  "

$ws.Range("E51").Value = "  +0.11%  "

# Restore default styling on the touched range so no stray style index is left behind.
$ws.Range("D2:E51").Style = "Normal"
